$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total_registros values for rows 2-6 (names unchanged)
$ws.Range("B2").Value = 142
$ws.Range("B3").Value = 142
$ws.Range("B4").Value = 139
$ws.Range("B5").Value = 127
$ws.Range("B6").Value = 120

# Rows 7 and 8 swap the empadronador names and get new totals
$ws.Range("A7").Value = "BURGA MEDINA SHIRLEY ROCIO"
$ws.Range("B7").Value = 114
$ws.Range("A8").Value = "BLANCO LOZANO ANDREA MILAGROS"
$ws.Range("B8").Value = 111

# Row 9 keeps its name, total updated
$ws.Range("B9").Value = 111
